# The commit adds the text "saddsasdaasd" to cell A1 of the first
# worksheet ("Plan1" / xl/worksheets/sheet1.xml), which is also the
# workbook's active sheet. Writing this text value causes the cell to be
# stored as a shared string and creates xl/sharedStrings.xml, matching
# the target workbook produced by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")
$ws.Range("A1").Value = "saddsasdaasd"
